$p = $ppt.ActivePresentation

# Slide 1: update the date in the subtitle placeholder
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(5).TextFrame.TextRange
$tr1.Runs(1, 1).Text = "Enterprise Solutions Architect | November 18, 2025"

# Slide 5: Implementation Approach - strip the markdown-style asterisks
# from the " *(Weeks n-n)*" run that follows each bold phase heading.
$s5 = $p.Slides.Item(5)
$tr = $s5.Shapes.Item(3).TextFrame.TextRange

$tr.Paragraphs(1, 1).Runs(2, 1).Text = " (Weeks 1-3)"
$tr.Paragraphs(6, 1).Runs(2, 1).Text = " (Weeks 4-6)"
$tr.Paragraphs(11, 1).Runs(2, 1).Text = " (Weeks 7-10)"
